# Updated cryptos list with refreshed price/volume data and a row 34/35
# swap (WEMIXToken now ranks above BinanceUSD, both with new figures).
#
# Note: column D ("Price") holds numeric-looking values (e.g. "247.10",
# "1.899.28") that must stay plain text, matching the source data feed's
# formatting (some have multiple "." separators, others would lose a
# trailing zero if Excel auto-converted them to a real number). Prefixing
# the literal with a leading single-quote forces Excel to store it as text
# without the quote character itself becoming part of the value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 34/35 swap: WEMIXToken now ranks above BinanceUSD, with fresh values.
$ws.Range("B34").Value = "WEMIXToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D34").Value = "'1.85"
$ws.Range("E34").Value = "  +14.83%  "

$ws.Range("B35").Value = "BinanceUSD"
$ws.Range("C35").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = "  +0.01%  "

# Price (column D) and Volume(1h) (column E) refreshes for the rest of the rows.
$ws.Range("D2").Value = "'35.651.57"
$ws.Range("E2").Value = "  +0.50%  "
$ws.Range("D3").Value = "'1.899.28"
$ws.Range("E3").Value = "  +0.19%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'247.10"
$ws.Range("E5").Value = "  -0.24%  "
$ws.Range("D6").Value = "'0.693"
$ws.Range("E6").Value = "  +0.17%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "'43.13"
$ws.Range("E8").Value = "  -2.02%  "
$ws.Range("D9").Value = "'57.28"
$ws.Range("E9").Value = "  +10.17%  "
$ws.Range("D10").Value = "'0.358"
$ws.Range("E10").Value = "  +1.76%  "
$ws.Range("D11").Value = "'0.0756"
$ws.Range("E11").Value = "  +2.24%  "
$ws.Range("E12").Value = "  +1.57%  "
$ws.Range("D13").Value = "'14.50"
$ws.Range("E13").Value = "  +10.79%  "
$ws.Range("D14").Value = "'0.805"
$ws.Range("E14").Value = "  +11.57%  "
$ws.Range("D15").Value = "'2.176.54"
$ws.Range("E15").Value = "  +0.14%  "
$ws.Range("D16").Value = "'5.04"
$ws.Range("E16").Value = "  +2.65%  "
$ws.Range("D17").Value = "'1.892.07"
$ws.Range("E17").Value = "  +0.27%  "
$ws.Range("D18").Value = "'35.657.71"
$ws.Range("E18").Value = "  +0.48%  "
$ws.Range("D19").Value = "'73.88"
$ws.Range("E19").Value = "  +0.45%  "
$ws.Range("E20").Value = "  +1.64%  "
$ws.Range("D21").Value = "'246.82"
$ws.Range("E21").Value = "  -0.36%  "
$ws.Range("D22").Value = "'13.01"
$ws.Range("E22").Value = "  +1.57%  "
$ws.Range("D23").Value = "'5.20"
$ws.Range("E23").Value = "  +5.01%  "
$ws.Range("E24").Value = "  +5.67%  "
$ws.Range("E25").Value = "  +0.03%  "
$ws.Range("E26").Value = "  -1.33%  "
$ws.Range("D27").Value = "'166.59"
$ws.Range("E27").Value = "  +0.61%  "
$ws.Range("D28").Value = "'8.65"
$ws.Range("E28").Value = "  +2.34%  "
$ws.Range("D29").Value = "'18.41"
$ws.Range("E29").Value = "  +0.34%  "
$ws.Range("E30").Value = "  +0.53%  "
$ws.Range("D31").Value = "'4.38"
$ws.Range("E31").Value = "  +3.19%  "
$ws.Range("E32").Value = "  +4.39%  "
$ws.Range("D33").Value = "'4.26"
$ws.Range("E36").Value = "  -16.57%  "
$ws.Range("D37").Value = "'0.857"
$ws.Range("E37").Value = "  +0.37%  "
$ws.Range("E38").Value = "  -1.10%  "
$ws.Range("D39").Value = "'0.0735"
$ws.Range("E39").Value = "  +8.09%  "
$ws.Range("D40").Value = "'0.0226"
$ws.Range("E40").Value = "  +6.41%  "
$ws.Range("D41").Value = "'99.37"
$ws.Range("E41").Value = "  +2.06%  "
$ws.Range("D42").Value = "'17.13"
$ws.Range("E42").Value = "  +0.14%  "
$ws.Range("D43").Value = "'14.62"
$ws.Range("E43").Value = "  +20.88%  "
$ws.Range("E44").Value = "  -0.03%  "
$ws.Range("D45").Value = "'1.323.27"
$ws.Range("E45").Value = "  +2.84%  "
$ws.Range("D46").Value = "'2.38"
$ws.Range("E46").Value = "  +1.13%  "
$ws.Range("D47").Value = "'0.0813"
$ws.Range("E47").Value = "  +1.64%  "
$ws.Range("E48").Value = "  -0.31%  "
$ws.Range("E49").Value = "  +0.12%  "
$ws.Range("E50").Value = "  +1.10%  "
$ws.Range("D51").Value = "'42.71"
$ws.Range("E51").Value = "  -0.97%  "
